$d = $word.ActiveDocument

# --- Step 1: Update the first paragraph's text ----------------------------
$r = $d.Range(0, 8)
$r.Text = "Hola mi nombre es Micaelo"

# --- Step 2: Drop a temporary bookmark right after "...Micaelo" -----------
# Placing a bookmark at this text position before any further edits stops
# the engine from coalescing the upcoming new run with this one when the
# document is serialized, so "Hola mi nombre es Micaelo" and the following
# " " stay as two distinct runs (matching the target OOXML).
$p1 = $d.Paragraphs(1)
$splitPos1 = $p1.Range.End - 1
[void]$d.Bookmarks.Add("TempSplit", $d.Range($splitPos1, $splitPos1))

# --- Step 3: Merge paragraph 1 ("...Micaelo") with paragraph 2 ("Yo") -----
# by deleting the paragraph mark that ends paragraph 1.
$mark1 = $d.Range($p1.Range.End - 1, $p1.Range.End)
[void]$mark1.Delete()

# --- Step 4: Merge the resulting paragraph with the old paragraph 3 -------
# ("Daniel") by deleting the paragraph mark that now ends paragraph 1.
$p1b = $d.Paragraphs(1)
$mark2 = $d.Range($p1b.Range.End - 1, $p1b.Range.End)
[void]$mark2.Delete()

# --- Step 5: Replace the old "Yo" text with a single space ----------------
$yoRange = $d.Content
[void]$yoRange.Find.Execute("Yo")
$yoStart = $yoRange.Start
$yoEnd = $yoRange.End
$d.Range($yoStart, $yoEnd).Text = " "

# The temporary bookmark has done its job; remove it again.
$temp = $d.Bookmarks("TempSplit")
[void]$temp.Delete()

# --- Step 6: Locate "Daniel" (where the new trailing sentence must start) -
$danRange = $d.Content
[void]$danRange.Find.Execute("Daniel")
$danStart = $danRange.Start
$danEnd = $danRange.End

# --- Step 7: Move the _GoBack bookmark to that split point first ----------
# (same trick as step 2: keeps the " " run and the new trailing-sentence
# run from being coalesced into one run).
$old = $d.Bookmarks("_GoBack")
[void]$old.Delete()
[void]$d.Bookmarks.Add("_GoBack", $d.Range($danStart, $danStart))

# --- Step 8: Replace "Daniel" with the new trailing sentence --------------
$trailer = ", en este documento estaré realizando unas cuantas pruebas para corroborar el funcionamiento de este proyecto . Pondré una palabra repetida , así como documento ."
$d.Range($danStart, $danEnd).Text = $trailer
